# Daily attendance processing - 2026-01-10 03:24:57
#
# The "Recorded By" column (G) lists the user(s) who recorded/updated a
# session, separated by ", ". The automated "System" account name was
# sometimes listed ahead of the real user's e-mail address. Normalize the
# ordering so that any "System" / "system" token is moved to the end of
# the list (capitalized "System" sorted ahead of lowercase "system" when
# both appear), while the relative order of the other entries (real users)
# is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("G$row")
    $value = $cell.Value2

    if (!$value) { continue }

    $parts = $value -split ", "
    if ($parts.Count -lt 2) { continue }

    $others = @()
    $systemParts = @()
    foreach ($part in $parts) {
        if ($part.ToLower() -eq "system") {
            $systemParts += $part
        } else {
            $others += $part
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    # Put capitalized "System" ahead of lowercase "system" within the
    # trailing group of system tokens.
    $systemUpper = @()
    $systemLower = @()
    foreach ($s in $systemParts) {
        if ($s.Equals("System")) {
            $systemUpper += $s
        } else {
            $systemLower += $s
        }
    }

    $newParts = $others + $systemUpper + $systemLower
    $newValue = $newParts -join ", "

    if (!$newValue.Equals($value)) {
        $cell.Value = $newValue
    }
}
